$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53:56 down to 54:57
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly data
$ws.Range("A53").Value = 6
$ws.Range("B53").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C53").Value = "Metropolitana"
$ws.Range("D53").Value = 44516
$ws.Range("D53").NumberFormat = $ws.Range("D54").NumberFormat
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = 100114007
$ws.Range("G53").Value = "Jengibre"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 400
$ws.Range("K53").Value = 13000
$ws.Range("L53").Value = 15000
$ws.Range("M53").Value = 14150
$ws.Range("N53").Value = "$/caja 13 kilos"
$ws.Range("O53").Value = "Perú"
$ws.Range("P53").Value = 1088
$ws.Range("Q53").Value = 13
$ws.Range("R53").Value = "Hortaliza"
